$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 650
$ws.Range("I49").Value = 550
$ws.Range("K49").Value = 1650
$ws.Range("M49").Value = -1514
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H86").Value = 2508.2
$ws.Range("I86").Value = 2013.2
$ws.Range("J86").Value = 3003.2
$ws.Range("K86").Value = 2013.2
$ws.Range("L86").Value = 3003.2
$ws.Range("M86").Value = -890.2
$ws.Range("N86").Value = -5249.2
$ws.Range("H89").Value = 2508.2
$ws.Range("I89").Value = 2013.2
$ws.Range("J89").Value = 3003.2
$ws.Range("K89").Value = 10066
$ws.Range("L89").Value = 15016
$ws.Range("M89").Value = -4450
$ws.Range("N89").Value = -26248
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H131").Value = 1207.6136
$ws.Range("I131").Value = 880.8823
$ws.Range("J131").Value = 1413.3334
$ws.Range("K131").Value = 2642.6469
$ws.Range("L131").Value = 4240.0002
$ws.Range("M131").Value = 2397.3531
$ws.Range("N131").Value = -14320.0002
$ws.Range("H138").Value = 3776747.2
$ws.Range("J138").Value = 4084351
$ws.Range("L138").Value = 12253053
$ws.Range("N138").Value = -12263333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 666.6667
$ws.Range("I3").Value = 666.6667
$ws.Range("K3").Value = 666.6667
$ws.Range("M3").Value = -551.6667
$ws.Range("H29").Value = 35000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H32").Value = 25121.846
$ws.Range("I32").Value = 24603.223
$ws.Range("J32").Value = 26288.75
$ws.Range("K32").Value = 24603.223
$ws.Range("L32").Value = 26288.75
$ws.Range("M32").Value = -24316.223
$ws.Range("N32").Value = -26862.75
$ws.Range("H74").Value = 5410284.5
$ws.Range("I74").Value = 7599628
$ws.Range("K74").Value = 7599628
$ws.Range("M74").Value = -7598754
$ws.Range("H77").Value = 5410284.5
$ws.Range("I77").Value = 7599628
$ws.Range("K77").Value = 37998140
$ws.Range("M77").Value = -37993772
$ws.Range("H122").Value = 3902.4
$ws.Range("J122").Value = 4125
$ws.Range("L122").Value = 12375
$ws.Range("N122").Value = -17275
$ws.Range("H132").Value = 35839.418
$ws.Range("I132").Value = 22803.25
$ws.Range("J132").Value = 80534.86
$ws.Range("K132").Value = 68409.75
$ws.Range("L132").Value = 241604.58
$ws.Range("M132").Value = -65879.75
$ws.Range("N132").Value = -246664.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 14930
$ws.Range("I7").Value = 2900
$ws.Range("J7").Value = 38990
$ws.Range("K7").Value = 2900
$ws.Range("L7").Value = 38990
$ws.Range("M7").Value = -2787
$ws.Range("N7").Value = -39216
$ws.Range("H29").Value = 703.2
$ws.Range("I29").Value = 703.2
$ws.Range("K29").Value = 703.2
$ws.Range("M29").Value = -414.2
$ws.Range("H122").Value = 53958.332
$ws.Range("J122").Value = 53958.332
$ws.Range("L122").Value = 53958.332
$ws.Range("N122").Value = -63758.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 243730.3
$ws.Range("I31").Value = 34627.168
$ws.Range("J31").Value = 557385
$ws.Range("K31").Value = 34627.168
$ws.Range("L31").Value = 557385
$ws.Range("M31").Value = -34332.168
$ws.Range("N31").Value = -557975
$ws.Range("H34").Value = 243730.3
$ws.Range("I34").Value = 34627.168
$ws.Range("J34").Value = 557385
$ws.Range("K34").Value = 34627.168
$ws.Range("L34").Value = 557385
$ws.Range("M34").Value = -34425.168
$ws.Range("N34").Value = -557789
$ws.Range("H58").Value = 24391662
$ws.Range("I58").Value = 30304236
$ws.Range("J58").Value = 2294.5
$ws.Range("K58").Value = 30304236
$ws.Range("L58").Value = 2294.5
$ws.Range("M58").Value = -30304033
$ws.Range("N58").Value = -2700.5
$ws.Range("H94").Value = 4786.5625
$ws.Range("I94").Value = 9087.857
$ws.Range("J94").Value = 1441.1111
$ws.Range("K94").Value = 9087.857
$ws.Range("L94").Value = 1441.1111
$ws.Range("M94").Value = -8636.857
$ws.Range("N94").Value = -2343.1111
$ws.Range("H107").Value = 417.97958
$ws.Range("I107").Value = 375.67648
$ws.Range("K107").Value = 375.67648
$ws.Range("M107").Value = 1544.32352
$ws.Range("H134").Value = 44157.72
$ws.Range("I134").Value = 854.5714
$ws.Range("J134").Value = 60997.832
$ws.Range("K134").Value = 2563.7142
$ws.Range("L134").Value = 182993.496
$ws.Range("M134").Value = -28.71420000000035
$ws.Range("N134").Value = -188063.496
$ws.Range("H136").Value = 24391662
$ws.Range("I136").Value = 30304236
$ws.Range("J136").Value = 2294.5
$ws.Range("K136").Value = 90912708
$ws.Range("L136").Value = 6883.5
$ws.Range("M136").Value = -90910158
$ws.Range("N136").Value = -11983.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 786.0323
$ws.Range("I131").Value = 295.54544
$ws.Range("J131").Value = 1055.8
$ws.Range("K131").Value = 886.63632
$ws.Range("L131").Value = 3167.4
$ws.Range("M131").Value = 4153.36368
$ws.Range("N131").Value = -13247.4
$ws.Range("H132").Value = 2376.3044
$ws.Range("I132").Value = 1480.3334
$ws.Range("J132").Value = 4056.25
$ws.Range("K132").Value = 13323.0006
$ws.Range("L132").Value = 36506.25
$ws.Range("M132").Value = -10793.0006
$ws.Range("N132").Value = -41566.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H122").Value = 2578.077
$ws.Range("I122").Value = 1950.875
$ws.Range("J122").Value = 3581.6
$ws.Range("K122").Value = 5852.625
$ws.Range("L122").Value = 10744.8
$ws.Range("M122").Value = -3402.625
$ws.Range("N122").Value = -15644.8
$ws.Range("H126").Value = 2108.6365
$ws.Range("I126").Value = 1562.2
$ws.Range("J126").Value = 2564
$ws.Range("K126").Value = 4686.6
$ws.Range("L126").Value = 7692
$ws.Range("M126").Value = -2216.6
$ws.Range("N126").Value = -12632

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2513.5334
$ws.Range("I40").Value = 2558.5833
$ws.Range("J40").Value = 2333.3333
$ws.Range("K40").Value = 2558.5833
$ws.Range("L40").Value = 2333.3333
$ws.Range("M40").Value = -2422.5833
$ws.Range("N40").Value = -2605.3333
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 19256.033
$ws.Range("I132").Value = 2002.2174
$ws.Range("J132").Value = 80308
$ws.Range("K132").Value = 6006.6522
$ws.Range("L132").Value = 240924
$ws.Range("M132").Value = -3476.6522
$ws.Range("N132").Value = -245984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1550001.5
$ws.Range("I2").Value = 3000000
$ws.Range("K2").Value = 3000000
$ws.Range("M2").Value = -2999888
$ws.Range("H132").Value = 38027.055
$ws.Range("I132").Value = 22943.623
$ws.Range("J132").Value = 113444.22
$ws.Range("K132").Value = 68830.86900000001
$ws.Range("L132").Value = 340332.66
$ws.Range("M132").Value = -66300.86900000001
$ws.Range("N132").Value = -345392.66
$ws.Range("H136").Value = 51212.855
$ws.Range("I136").Value = 38322.85
$ws.Range("J136").Value = 76072.14
$ws.Range("K136").Value = 114968.55
$ws.Range("L136").Value = 228216.42
$ws.Range("M136").Value = -112418.55
$ws.Range("N136").Value = -233316.42
$ws.Range("H139").Value = 21666.666
$ws.Range("I139").Value = 15000
$ws.Range("K139").Value = 15000
$ws.Range("M139").Value = -9860
$ws.Range("H141").Value = 65871.58
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 67864.44500000001
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 67864.44500000001
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -78224.44500000001
